$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new values look like plain numbers (e.g. "1.003",
# "152.50") must be forced to Text format before the value is written,
# otherwise Excel auto-converts them to numbers and formatting like
# trailing zeros would be lost (e.g. "152.50" -> 152.5). Cells whose new
# value already cannot parse as a number (e.g. "28.584.60") are left alone.
$textCells = @("D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D19","D22","D23","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.584.60'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').Value = '1.867.92'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '324.12'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4606'
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('D8').Value = '0.3871'
$ws.Range('D9').Value = '0.07855'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').Value = '0.9737'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').Value = '21.89'
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').Value = '1.871.69'
$ws.Range('E12').Value = '  +4.12%  '
$ws.Range('D13').Value = '6.975'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').Value = '5.693'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').Value = '0.06947'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').Value = '88.12'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').Value = '16.78'
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').Value = '28.591.90'
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('D22').Value = '5.268'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '11.01'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').Value = '2.108.95'
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('D26').Value = '152.50'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = '19.19'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').Value = '5.768'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').Value = '1.985'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('D30').Value = '119.07'
$ws.Range('E30').Value = '  +1.43%  '
$ws.Range('D31').Value = '0.09334'
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').Value = '0.9188'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('D33').Value = '5.260'
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').Value = '1.333'
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('D35').Value = '3.325'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').Value = '0.05790'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('D37').Value = '1.155'
$ws.Range('E37').Value = '  +1.18%  '
$ws.Range('D38').Value = '0.02084'
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').Value = '7.737'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').Value = '0.5624'
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').Value = '0.1783'
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('D42').Value = '9.777'
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D43').Value = '0.07179'
$ws.Range('E43').Value = '  +2.45%  '
$ws.Range('D44').Value = '11.69'
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('D45').Value = '0.5287'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = '2.147'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').Value = '1.135'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').Value = '1.832'
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('D49').Value = '112.79'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').Value = '2.410'
$ws.Range('E50').Value = '  +3.80%  '
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.15%  '
